# issue #5: stock data from json to db
# Update the "股票" (stock) worksheet (sheet 5): insert a "category" column
# after property_category, a "source_file" and "index" column after
# legislator_id, fix the malformed face value for 聯華電子, and keep it as
# text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# ---- Header row (row 1) ----
# Existing headers B1:K1 (name, owner, quantity, face_value, currency,
# total, property_category, date, legislator_name, legislator_id) stay put;
# shift the trailing two (date, legislator_name, legislator_id) out to make
# room for the new "category" column, then append the two new trailing
# columns.
$ws.Cells.Item(1, 9).Value  = "category"        # I1
$ws.Cells.Item(1, 10).Value = "date"             # J1
$ws.Cells.Item(1, 11).Value = "legislator_name"  # K1
$ws.Cells.Item(1, 12).Value = "legislator_id"    # L1
$ws.Cells.Item(1, 13).Value = "source_file"      # M1
$ws.Cells.Item(1, 14).Value = "index"            # N1

# ---- Data rows ----
$rows = 2, 3, 4, 5
foreach ($r in $rows) {
    # H (property_category) keeps "stock"; insert "normal" as the new
    # category value in I, then shift date/legislator_name/legislator_id
    # out to J/K/L and append source_file/index in M/N.
    $ws.Cells.Item($r, 8).Value  = "stock"
    $ws.Cells.Item($r, 9).Value  = "normal"
    # Force text so "2011-11-21" is not auto-converted to a date serial.
    $ws.Cells.Item($r, 10).NumberFormat = "@"
    $ws.Cells.Item($r, 10).Value = "2011-11-21"
    $ws.Cells.Item($r, 11).Value = "蔣乃辛"
    $ws.Cells.Item($r, 12).Value = 1722
    $ws.Cells.Item($r, 13).Value = "tmp12421"
}

# index column (N) mirrors the existing row index in column A
$ws.Cells.Item(2, 14).Value = 58
$ws.Cells.Item(3, 14).Value = 59
$ws.Cells.Item(4, 14).Value = 60
$ws.Cells.Item(5, 14).Value = 61

# Fix 聯華電子's face value: "279；030" (with a stray full-width
# semicolon) -> "279030", keeping it stored as text.
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "279030"
